# Add Denmark, Sweden and Norway market test-data sheets, cloned from the
# existing "Belgium" sheet (same layout/styles/merges), with the market name
# and NGC code swapped in B2 / B4 respectively. Norway ends up the active /
# selected sheet (last one added), matching the workbook-level activeTab.

$wb = $excel.ActiveWorkbook
$belgium = $wb.Worksheets.Item("Belgium")

# --- Denmark -------------------------------------------------------------
$belgium.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$denmark = $wb.Worksheets.Item($wb.Worksheets.Count)
$denmark.Name = "Denmark"
$denmark.Range("B2").Value = "Denmark Market"
$denmark.Range("B4").Value = "NGC-3446/T2004"
$denmark.Range("A1:XFD1048576").Select()

# --- Sweden ----------------------------------------------------------------
$belgium.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sweden = $wb.Worksheets.Item($wb.Worksheets.Count)
$sweden.Name = "Sweden"
$sweden.Range("B2").Value = "Sweden Market"
$sweden.Range("B4").Value = "NGC-3465/T2025"
$sweden.Range("A1:XFD1048576").Select()

# --- Norway (left active/selected, like in the authored workbook) --------
$belgium.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$norway = $wb.Worksheets.Item($wb.Worksheets.Count)
$norway.Name = "Norway"
$norway.Range("B2").Value = "Norway Market"
$norway.Range("B4").Value = "NGC-3464/T1919"
$norway.Activate()
$norway.Range("B2:B4").Select()
